$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'318.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'4.38%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'36.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.21%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.132"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.22%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'4.43%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.154"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.29%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.009"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "'0.89%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9264"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.63%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1005"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.76%"
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'1.32%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09225"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'6.77%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03606"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'3.39%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09931"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.02%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001442"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.67%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005687"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.19%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.464"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.09%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'13.59%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3373"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.54%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1302"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-2.00%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.060"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'4.48%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2190"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.46%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04595"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.14%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'1.14%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004732"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-7.05%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001250"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-10.62%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0004505"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-5.17%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02011"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'8.84%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04988"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'4.19%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007804"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.30%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1402"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.25%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007816"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'1.21%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002097"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-5.91%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01182"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'7.38%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006470"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'1.17%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.08%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'18.20%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.001902"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-4.91%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.08%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.08%"
$ws.Range("E51").Style = "Normal"
